# Cover Letter edit: update the three body paragraphs (interest, experience,
# and communication/closing paragraphs) to match the revised wording, and
# simplify the final "Thank you" paragraph into a single run.

$d = $word.ActiveDocument

function Set-ParagraphText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    # Exclude the trailing paragraph mark so we don't delete the paragraph itself.
    $textRange = $d.Range($r.Start, $r.End - 1)
    # Clear the paragraph first, then insert the new text. This guarantees the
    # paragraph collapses to a single run instead of retaining old run breaks
    # when substrings happen to match the replacement text.
    $textRange.Text = ""
    $p2 = $d.Paragraphs($paraIndex)
    $insertStart = $p2.Range.Start
    $insertRange = $d.Range($insertStart, $insertStart)
    $insertRange.Text = $newText
}

# Paragraph 3: "I am writing to express my interest ..."
Set-ParagraphText 3 "I am writing to express my interest in the position at your company. As a recent graduate of Bradley University with a degree in data science, and with over a year and a half of experience in the industry, I am confident that my skills and experience align well with the requirements of this role."

# Paragraph 5: "During my time at Produce Pro Software, ..." -> "In my previous positions, ..."
Set-ParagraphText 5 "In my previous positions, I gained extensive experience in both data analytics and administration. My work in the analytics department at Produce Pro involved maintaining data warehouses for over 60 customers and supporting their use of Cognos, a business intelligence tool. After about a year I was brought on by our parent company to do a similar job but using their Windows/PowerBI solutions and had to learn a new set of tools for that position. These different roles allowed me to work with multiple customers with multiple analytics products and gain a deep understanding of both how data can be used to drive business decisions and how diverse companies needs for that data can be. "

# Paragraph 7: "In addition to my technical skills, ..."
$rsquo = [char]0x2019
Set-ParagraphText 7 "In addition to my technical skills, I am a strong communicator and collaborator. In both of my positions I worked closely with both technical and non-technical stakeholders to understand their needs and deliver solutions that met their requirements. I${rsquo}m also a fast learner and am always eager to take on new challenges and expand my skills. I believe that I can bring this experience and energy to your organization to help build up your analytics process and drive insights for your company."

# Paragraph 9: "Thank you for considering my application. ..." -> merge into a single run
Set-ParagraphText 9 "Thank you for considering my application. I would love to learn more about the opportunities available at your company and discuss my qualifications further."

Write-Output "Done"
